# Data.xlsx upload edit:
#  - Fix the Tamil transliteration text for "Gnana Nirmalame" (TA sheet, D2):
#       ஞான நிர்மலாமே  ->  ஞான நிர்மலமே
#  - Resize a couple of "best fit" columns (HI!C and TA!C:D) to their new
#    auto-fit widths.
#  - Leave the cursor/selection on MR, HI and (finally, as the active tab) TA
#    at the cell positions the author ended up on.

$wb = $excel.ActiveWorkbook

# ---- TA sheet: correct the Tamil spelling in D2 --------------------------
$wsTA = $wb.Worksheets.Item("TA")
$wsTA.Range("D2").Value2 = "ஞான நிர்மலமே"

# ---- HI sheet: widen column C to its new best-fit width ------------------
$wsHI = $wb.Worksheets.Item("HI")
$wsHI.Columns.Item(3).ColumnWidth = 17.5

# ---- TA sheet: give columns C and D explicit best-fit widths -------------
$wsTA.Columns.Item(3).ColumnWidth = 15.17
$wsTA.Columns.Item(4).ColumnWidth = 19.33

# ---- Selections -------------------------------------------------------
# MR sheet: move the selection to C18 (leave MR as a non-active tab)
$wsMR = $wb.Worksheets.Item("MR")
$wsMR.Range("C18").Select()

# HI sheet: move the selection to D9 (leave HI as a non-active tab)
$wsHI.Range("D9").Select()

# TA sheet: move the selection to D3; TA stays the active/selected tab
$wsTA.Range("D3").Select()
